$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.649.17"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.449.17"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'578.46"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'148.54"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "'8.05"
$ws.Range("E9").Value = "  +5.31%  "
$ws.Range("D11").Value = "'0.413"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").Value = "4.040.14"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "'28.31"
$ws.Range("E14").Value = "  -5.54%  "
$ws.Range("D15").Value = "3.450.39"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "62.717.83"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'14.61"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "'8.99"
$ws.Range("E20").Value = "  -3.09%  "
$ws.Range("D21").Value = "'387.06"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'75.31"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "3.584.01"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "'0.182"
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'7.98"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -4.28%  "
$ws.Range("D34").Value = "'23.23"
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").Value = "'5.38"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").Value = "'1.63"
$ws.Range("E36").Value = "  +3.49%  "
$ws.Range("D37").Value = "'32.13"
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "'169.20"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "3.484.15"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").Value = "'0.0777"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "'0.785"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "'4.38"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "2.569.18"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "'6.90"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "'2.24"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'22.53"
$ws.Range("E50").Value = "  -4.15%  "
$ws.Range("E51").Value = "  +0.00%  "
